$wb = $excel.ActiveWorkbook

# --- Login sheet: selection changes from A1:B2 range-select to a single-cell
#     select on B2 (activeCell=B2, sqref=B2) ---
$wsLogin = $wb.Worksheets.Item("Login")
$wsLogin.Range("B2").Select()

# --- Register sheet: update the credential data (email addresses bumped
#     from "002" to "003"), move the selection to D15, and make this the
#     active/visible tab ---
$wsRegister = $wb.Worksheets.Item("Register")
$wsRegister.Range("D2").Value = "david003@manchesterunited.com"
$wsRegister.Range("D3").Value = "victoriabeckham003@spicegirls.com"
$wsRegister.Range("D15").Select()
$wsRegister.Activate()

# --- AddToCart sheet: password cell now reuses the existing "Password@123"
#     shared string instead of the old "Passwor@123" typo'd one (which gets
#     dropped from the shared-string table entirely since it becomes unused) ---
$wsCart = $wb.Worksheets.Item("AddToCart")
$wsCart.Range("B2").Value = "Password@123"
